# OPT-2 / OPT-15: update field names/labels across reference sheets and
# introduce several new columns (database fields) required by the fixed
# target-function calculators / excel parsers / templates.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "FilmRecipes" (index 2): rename a few headers (unit switch to the
# euro sign and DN instead of mk).
# ---------------------------------------------------------------------
$wsFilmRecipes = $wb.Worksheets.Item(2)
$wsFilmRecipes.Range("F1").Value = "Nozzle, DN"
$wsFilmRecipes.Range("G1").Value = "Calibration, DN"
$wsFilmRecipes.Range("H1").Value = "Cooling lip, DN"

# ---------------------------------------------------------------------
# Sheet "Customers" (index 3): add a "Number" column.
# ---------------------------------------------------------------------
$wsCustomers = $wb.Worksheets.Item(3)
$wsCustomers.Range("B1").Value = "Number"
$wsCustomers.Range("A1").Copy()
$wsCustomers.Range("B1").PasteSpecial(-4122)
$wsCustomers.Columns.Item(2).ColumnWidth = 15.333333333333334

# ---------------------------------------------------------------------
# Sheet "Orders" (index 4): rename / insert new production-tracking
# columns (FinishedGoods, Waste, RollsCount, PredefinedTime) and move the
# remaining columns after them.
# ---------------------------------------------------------------------
$wsOrders = $wb.Worksheets.Item(4)
$wsOrders.Range("B1").Value = "Width"
$wsOrders.Range("D1").Value = "FinishedGoods, kg"
$wsOrders.Range("E1").Value = "Waste, kg"
$wsOrders.Range("F1").Value = "RollsCount"
$wsOrders.Range("G1").Value = "PredefinedTime, min"
$wsOrders.Range("H1").Value = "Film recipe name"
$wsOrders.Range("I1").Value = "Planning end date"
$wsOrders.Range("J1").Value = "Price overdue, EUR/h"
$wsOrders.Range("K1").Value = "Customer name"
$wsOrders.Range("D1").Copy()
$wsOrders.Range("J1:K1").PasteSpecial(-4122)
$wsOrders.Columns.Item(10).ColumnWidth = 17.333333333333332
$wsOrders.Columns.Item(11).ColumnWidth = 13.666666666666666

# Fix the currency sign on the price overdue label (€ instead of EUR).
$wsOrders.Range("J1").Value = [char]0x20AC + "/h"
$wsOrders.Range("J1").Value = "Price overdue, " + [char]0x20AC + "/h"

# ---------------------------------------------------------------------
# Sheet "Extruders" (index 5): rename width fields to mm, and add new
# consumption columns for thickness/width change.
# ---------------------------------------------------------------------
$wsExtruders = $wb.Worksheets.Item(5)
$wsExtruders.Range("E1").Value = "Width min, mm"
$wsExtruders.Range("F1").Value = "Width max, mm"
$wsExtruders.Range("M1").Value = "Thickness change time, min"
$wsExtruders.Range("N1").Value = "Thickness change consumption, kg/h"
$wsExtruders.Range("O1").Value = "Width change time, min"
$wsExtruders.Range("P1").Value = "Width change consumption, kg/h"
$wsExtruders.Range("M1").Copy()
$wsExtruders.Range("N1:P1").PasteSpecial(-4122)
$wsExtruders.Columns.Item(14).ColumnWidth = 22.0
$wsExtruders.Columns.Item(15).ColumnWidth = 23.5
$wsExtruders.Columns.Item(16).ColumnWidth = 20.0

# ---------------------------------------------------------------------
# Sheet "FilmTypesChanges" (index 6): rename the "article" change columns
# to "recipe name" ones, and add a Consumption column.
# ---------------------------------------------------------------------
$wsFilmTypesChanges = $wb.Worksheets.Item(6)
$wsFilmTypesChanges.Range("B1").Value = "Name recipe from change"
$wsFilmTypesChanges.Range("C1").Value = "Name recipe to change"
$wsFilmTypesChanges.Range("E1").Value = "Consumption, kg/h"
$wsFilmTypesChanges.Range("A1").Copy()
$wsFilmTypesChanges.Range("E1").PasteSpecial(-4122)
$wsFilmTypesChanges.Columns.Item(5).ColumnWidth = 21.666666666666668

# ---------------------------------------------------------------------
# Sheets "NozzleChanges" / "CalibrationChanges" / "CoolingLipChanges"
# (indices 7, 8, 9): add a Consumption column to each.
# ---------------------------------------------------------------------
$wsNozzleChanges = $wb.Worksheets.Item(7)
$wsNozzleChanges.Range("D1").Value = "Consumption, kg/h"
$wsNozzleChanges.Range("A1").Copy()
$wsNozzleChanges.Range("D1").PasteSpecial(-4122)
$wsNozzleChanges.Columns.Item(4).ColumnWidth = 19.0

$wsCalibrationChanges = $wb.Worksheets.Item(8)
$wsCalibrationChanges.Range("D1").Value = "Consumption, kg/h"
$wsCalibrationChanges.Range("A1").Copy()
$wsCalibrationChanges.Range("D1").PasteSpecial(-4122)
$wsCalibrationChanges.Columns.Item(4).ColumnWidth = 20.666666666666668

$wsCoolingLipChanges = $wb.Worksheets.Item(9)
$wsCoolingLipChanges.Range("D1").Value = "Consumption, kg/h"
$wsCoolingLipChanges.Range("A1").Copy()
$wsCoolingLipChanges.Range("D1").PasteSpecial(-4122)
$wsCoolingLipChanges.Columns.Item(4).ColumnWidth = 20.5

# ---------------------------------------------------------------------
# Fix the Euro currency sign labels (material cost, price overdue),
# typed last so these two shared strings land at the end of the table,
# matching the editing order of the original author.
# ---------------------------------------------------------------------
$wsFilmRecipes.Range("E1").Value = "Material cost, " + [char]0x20AC + "/kg"

# Re-select cells roughly where the author left off on each sheet.
$wsFilmRecipes.Range("E1").Select()
$wsCustomers.Range("A1:B1").Select()
$wsOrders.Range("B28").Select()
$wsExtruders.Range("N9").Select()
$wsFilmTypesChanges.Range("A1:E1").Select()
$wsNozzleChanges.Range("D1").Select()
$wsCalibrationChanges.Range("D1").Select()
$wsCoolingLipChanges.Range("D1").Select()

$wb.Save()
